# Auto-generated edit script applying the Ultros_Profits.xlsx diff
# Updates numeric cells across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3500.6667
$ws.Range("I106").Value = 3500.6667
$ws.Range("K106").Value = 3500.6667
$ws.Range("M106").Value = -2869.6667

$ws.Range("H116").Value = 7500
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

$ws.Range("H125").Value = 2112.5
$ws.Range("I125").Value = 735.7778
$ws.Range("J125").Value = 3882.5715
$ws.Range("K125").Value = 6622.000199999999
$ws.Range("L125").Value = 34943.1435
$ws.Range("M125").Value = -4162.000199999999
$ws.Range("N125").Value = -39863.1435

$ws.Range("H138").Value = 2511.394
$ws.Range("J138").Value = 3024.7896
$ws.Range("L138").Value = 9074.3688
$ws.Range("N138").Value = -19354.3688

$ws.Range("H141").Value = 2666.4443
$ws.Range("I141").Value = 2725.8
$ws.Range("K141").Value = 8177.400000000001
$ws.Range("M141").Value = -2997.400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 18333.334
$ws.Range("I45").Value = 15000
$ws.Range("K45").Value = 15000
$ws.Range("M45").Value = -14623

$ws.Range("H61").Value = 3049.55
$ws.Range("I61").Value = 1745.7693
$ws.Range("J61").Value = 5470.857
$ws.Range("K61").Value = 1745.7693
$ws.Range("L61").Value = 5470.857
$ws.Range("M61").Value = -1533.7693
$ws.Range("N61").Value = -5894.857

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").ClearContents()
$ws.Range("N124").ClearContents()

$ws.Range("H128").Value = 60429
$ws.Range("J128").Value = 60429
$ws.Range("L128").Value = 60429
$ws.Range("N128").Value = -70389

$ws.Range("H132").Value = 4735.3696
$ws.Range("J132").Value = 4803.3335
$ws.Range("L132").Value = 14410.0005
$ws.Range("N132").Value = -19470.0005

$ws.Range("H136").Value = 3049.55
$ws.Range("I136").Value = 1745.7693
$ws.Range("J136").Value = 5470.857
$ws.Range("K136").Value = 5237.3079
$ws.Range("L136").Value = 16412.571
$ws.Range("M136").Value = -2687.3079
$ws.Range("N136").Value = -21512.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 11074
$ws.Range("I97").Value = 11074
$ws.Range("K97").Value = 11074
$ws.Range("M97").Value = -10083

$ws.Range("H122").Value = 49999.77
$ws.Range("J122").Value = 49999.77
$ws.Range("L122").Value = 49999.77
$ws.Range("N122").Value = -59799.77

$ws.Range("H126").Value = 49914.445
$ws.Range("J126").Value = 49914.445
$ws.Range("L126").Value = 49914.445
$ws.Range("N126").Value = -59794.445

$ws.Range("H132").Value = 70222.22
$ws.Range("J132").Value = 70222.22
$ws.Range("L132").Value = 70222.22
$ws.Range("N132").Value = -80342.22

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 671.75
$ws.Range("I105").Value = 633.8
$ws.Range("K105").Value = 633.8
$ws.Range("M105").Value = 1113.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 107895230
$ws.Range("I4").Value = 205000140
$ws.Range("J4").Value = 894.1111
$ws.Range("K4").Value = 615000420
$ws.Range("L4").Value = 2682.3333
$ws.Range("M4").Value = -615000308
$ws.Range("N4").Value = -2906.3333

$ws.Range("H7").Value = 381.52173
$ws.Range("I7").Value = 293.625
$ws.Range("J7").Value = 582.4286
$ws.Range("K7").Value = 880.875
$ws.Range("L7").Value = 1747.2858
$ws.Range("M7").Value = -768.875
$ws.Range("N7").Value = -1971.2858

$ws.Range("H23").Value = 3324.2727
$ws.Range("J23").Value = 2721.1667
$ws.Range("L23").Value = 8163.500100000001
$ws.Range("N23").Value = -8633.500100000001

$ws.Range("H58").Value = 7428.4287
$ws.Range("J58").Value = 7428.4287
$ws.Range("L58").Value = 22285.2861
$ws.Range("N58").Value = -22541.2861

$ws.Range("H113").Value = 1440.1904
$ws.Range("I113").Value = 1695.5
$ws.Range("K113").Value = 5086.5
$ws.Range("M113").Value = -2916.5

$ws.Range("H129").Value = 1912.7333
$ws.Range("I129").Value = 1081.8
$ws.Range("J129").Value = 2328.2
$ws.Range("K129").Value = 3245.4
$ws.Range("L129").Value = 6984.599999999999
$ws.Range("M129").Value = 1754.6
$ws.Range("N129").Value = -16984.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1450711.9
$ws.Range("I3").Value = 1666.3334
$ws.Range("J3").Value = 2537496
$ws.Range("K3").Value = 1666.3334
$ws.Range("L3").Value = 2537496
$ws.Range("M3").Value = -1550.3334
$ws.Range("N3").Value = -2537728

$ws.Range("H18").Value = 5000
$ws.Range("J18").Value = 5000
$ws.Range("L18").Value = 5000
$ws.Range("N18").Value = -5586

$ws.Range("H20").Value = 14001000
$ws.Range("J20").Value = 5000
$ws.Range("L20").Value = 5000
$ws.Range("N20").Value = -5490

$ws.Range("H80").Value = 19676846
$ws.Range("I80").Value = 88742.234
$ws.Range("K80").Value = 88742.234
$ws.Range("M80").Value = -87744.234

$ws.Range("H83").Value = 19676846
$ws.Range("I83").Value = 88742.234
$ws.Range("K83").Value = 443711.17
$ws.Range("M83").Value = -438719.17

$ws.Range("H133").Value = 70376.71000000001
$ws.Range("J133").Value = 70376.71000000001
$ws.Range("L133").Value = 70376.71000000001
$ws.Range("N133").Value = -80496.71000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 1499.3334
$ws.Range("J25").Value = 1499.3334
$ws.Range("L25").Value = 1499.3334
$ws.Range("N25").Value = -1959.3334

$ws.Range("H43").Value = 2509502.8
$ws.Range("I43").Value = 12670.333
$ws.Range("J43").Value = 10000000
$ws.Range("K43").Value = 12670.333
$ws.Range("L43").Value = 10000000
$ws.Range("M43").Value = -12477.333
$ws.Range("N43").Value = -10000386

$ws.Range("H82").Value = 1478.125
$ws.Range("I82").Value = 794
$ws.Range("J82").Value = 2357.7144
$ws.Range("K82").Value = 794
$ws.Range("L82").Value = 2357.7144
$ws.Range("M82").Value = -433
$ws.Range("N82").Value = -3079.7144

$ws.Range("H85").Value = 1478.125
$ws.Range("I85").Value = 794
$ws.Range("J85").Value = 2357.7144
$ws.Range("K85").Value = 794
$ws.Range("L85").Value = 2357.7144
$ws.Range("M85").Value = 454
$ws.Range("N85").Value = -4853.7144

$ws.Range("H136").Value = 3973.0732
$ws.Range("I136").Value = 2869.5
$ws.Range("K136").Value = 8608.5
$ws.Range("M136").Value = -6058.5

$ws.Range("H139").Value = 92272
$ws.Range("J139").Value = 94444.44500000001
$ws.Range("L139").Value = 94444.44500000001
$ws.Range("N139").Value = -104724.445

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 33582.832
$ws.Range("J94").Value = 33582.832
$ws.Range("L94").Value = 33582.832
$ws.Range("N94").Value = -35384.832

$ws.Range("H126").Value = 2604.5833
$ws.Range("I126").Value = 2482.7778
$ws.Range("K126").Value = 7448.3334
$ws.Range("M126").Value = -4978.3334

$ws.Range("H136").Value = 9436983
$ws.Range("J136").Value = 4717.4116
$ws.Range("L136").Value = 14152.2348
$ws.Range("N136").Value = -19252.2348
